# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# to the Tonberry_Profits workbook (updates currentAveragePrice* / LevePrice* / LeveProfit* columns).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(42, 8).Value = 198.2
$ws.Cells.Item(42, 9).Value = 120.5
$ws.Cells.Item(42, 11).Value = 361.5
$ws.Cells.Item(42, 13).Value = -131.5
$ws.Cells.Item(57, 8).Value = 32000
$ws.Cells.Item(57, 10).Value = 32000
$ws.Cells.Item(57, 12).Value = 96000
$ws.Cells.Item(57, 14).Value = -96998
$ws.Cells.Item(112, 8).Value = 1594.1052
$ws.Cells.Item(112, 10).Value = 1594.1052
$ws.Cells.Item(112, 12).Value = 4782.3156
$ws.Cells.Item(112, 14).Value = -6998.3156
$ws.Cells.Item(113, 8).Value = 31461.5
$ws.Cells.Item(113, 9).Value = 61002
$ws.Cells.Item(113, 10).Value = 1921
$ws.Cells.Item(113, 11).Value = 61002
$ws.Cells.Item(113, 12).Value = 1921
$ws.Cells.Item(113, 13).Value = -57748
$ws.Cells.Item(113, 14).Value = -8429
$ws.Cells.Item(132, 8).Value = 1333.2941
$ws.Cells.Item(132, 10).Value = 1553
$ws.Cells.Item(132, 12).Value = 4659
$ws.Cells.Item(132, 14).Value = -9719

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(24, 8).Value = 0
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 14).ClearContents()
$ws.Cells.Item(32, 8).Value = 6599.775
$ws.Cells.Item(32, 9).Value = 4863.758
$ws.Cells.Item(32, 11).Value = 4863.758
$ws.Cells.Item(32, 13).Value = -4576.758
$ws.Cells.Item(74, 8).Value = 1001.1892
$ws.Cells.Item(74, 9).Value = 548.4828
$ws.Cells.Item(74, 10).Value = 2642.25
$ws.Cells.Item(74, 11).Value = 548.4828
$ws.Cells.Item(74, 12).Value = 2642.25
$ws.Cells.Item(74, 13).Value = 325.5172
$ws.Cells.Item(74, 14).Value = -4390.25
$ws.Cells.Item(77, 8).Value = 1001.1892
$ws.Cells.Item(77, 9).Value = 548.4828
$ws.Cells.Item(77, 10).Value = 2642.25
$ws.Cells.Item(77, 11).Value = 2742.414
$ws.Cells.Item(77, 12).Value = 13211.25
$ws.Cells.Item(77, 13).Value = 1625.586
$ws.Cells.Item(77, 14).Value = -21947.25
$ws.Cells.Item(82, 8).Value = 88888.336
$ws.Cells.Item(82, 9).Value = 66666
$ws.Cells.Item(82, 11).Value = 66666
$ws.Cells.Item(82, 13).Value = -66305
$ws.Cells.Item(85, 8).Value = 88888.336
$ws.Cells.Item(85, 9).Value = 66666
$ws.Cells.Item(85, 11).Value = 66666
$ws.Cells.Item(85, 13).Value = -65418
$ws.Cells.Item(96, 8).Value = 100000
$ws.Cells.Item(96, 10).Value = 100000
$ws.Cells.Item(96, 12).Value = 100000
$ws.Cells.Item(96, 14).Value = -105492
$ws.Cells.Item(100, 8).Value = 0
$ws.Cells.Item(100, 10).Value = 0
$ws.Cells.Item(100, 12).Value = 0
$ws.Cells.Item(100, 14).ClearContents()
$ws.Cells.Item(110, 8).Value = 193.6
$ws.Cells.Item(110, 9).Value = 193.6
$ws.Cells.Item(110, 11).Value = 193.6
$ws.Cells.Item(110, 13).Value = 1851.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1403.5714
$ws.Cells.Item(94, 9).Value = 665
$ws.Cells.Item(94, 11).Value = 665
$ws.Cells.Item(94, 13).Value = -214

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 14).ClearContents()
$ws.Cells.Item(31, 8).Value = 1733.3636
$ws.Cells.Item(31, 9).Value = 1342.9
$ws.Cells.Item(31, 10).Value = 2058.75
$ws.Cells.Item(31, 11).Value = 1342.9
$ws.Cells.Item(31, 12).Value = 2058.75
$ws.Cells.Item(31, 13).Value = -1047.9
$ws.Cells.Item(31, 14).Value = -2648.75
$ws.Cells.Item(34, 8).Value = 1733.3636
$ws.Cells.Item(34, 9).Value = 1342.9
$ws.Cells.Item(34, 10).Value = 2058.75
$ws.Cells.Item(34, 11).Value = 1342.9
$ws.Cells.Item(34, 12).Value = 2058.75
$ws.Cells.Item(34, 13).Value = -1140.9
$ws.Cells.Item(34, 14).Value = -2462.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 63.090908
$ws.Cells.Item(12, 9).Value = 39.166668
$ws.Cells.Item(12, 10).Value = 91.8
$ws.Cells.Item(12, 11).Value = 117.500004
$ws.Cells.Item(12, 12).Value = 275.4
$ws.Cells.Item(12, 13).Value = 55.499996
$ws.Cells.Item(12, 14).Value = -621.4
$ws.Cells.Item(64, 8).Value = 3249.75
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 13).ClearContents()
$ws.Cells.Item(67, 8).Value = 3249.75
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 11).Value = 0
$ws.Cells.Item(67, 13).ClearContents()
$ws.Cells.Item(109, 8).Value = 5069.357
$ws.Cells.Item(109, 9).Value = 1333
$ws.Cells.Item(109, 10).Value = 6088.364
$ws.Cells.Item(109, 11).Value = 3999
$ws.Cells.Item(109, 12).Value = 18265.092
$ws.Cells.Item(109, 13).Value = -2959
$ws.Cells.Item(109, 14).Value = -20345.092
$ws.Cells.Item(113, 8).Value = 6154.263
$ws.Cells.Item(113, 10).Value = 937.05884
$ws.Cells.Item(113, 12).Value = 2811.17652
$ws.Cells.Item(113, 14).Value = -7151.17652
$ws.Cells.Item(129, 8).Value = 61580.582
$ws.Cells.Item(129, 9).Value = 999
$ws.Cells.Item(129, 10).Value = 67088
$ws.Cells.Item(129, 11).Value = 2997
$ws.Cells.Item(129, 12).Value = 201264
$ws.Cells.Item(129, 13).Value = 2003
$ws.Cells.Item(129, 14).Value = -211264
$ws.Cells.Item(131, 8).Value = 26228.965
$ws.Cells.Item(131, 10).Value = 30491.291
$ws.Cells.Item(131, 12).Value = 91473.87300000001
$ws.Cells.Item(131, 14).Value = -101553.873

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 35000
$ws.Cells.Item(15, 10).Value = 35000
$ws.Cells.Item(15, 12).Value = 35000
$ws.Cells.Item(15, 14).Value = -35576
$ws.Cells.Item(81, 8).Value = 35000
$ws.Cells.Item(81, 10).Value = 35000
$ws.Cells.Item(81, 12).Value = 35000
$ws.Cells.Item(81, 14).Value = -36996
$ws.Cells.Item(84, 8).Value = 35000
$ws.Cells.Item(84, 10).Value = 35000
$ws.Cells.Item(84, 12).Value = 105000
$ws.Cells.Item(84, 14).Value = -114984
$ws.Cells.Item(97, 8).Value = 973.86664
$ws.Cells.Item(97, 9).Value = 1010.1539
$ws.Cells.Item(97, 11).Value = 1010.1539
$ws.Cells.Item(97, 13).Value = -514.1539
$ws.Cells.Item(126, 8).Value = 3144433.8
$ws.Cells.Item(126, 9).Value = 5053123
$ws.Cells.Item(126, 10).Value = 145065
$ws.Cells.Item(126, 11).Value = 15159369
$ws.Cells.Item(126, 12).Value = 435195
$ws.Cells.Item(126, 13).Value = -15156899
$ws.Cells.Item(126, 14).Value = -440135
$ws.Cells.Item(132, 8).Value = 5498119
$ws.Cells.Item(132, 9).Value = 12823061
$ws.Cells.Item(132, 10).Value = 4412
$ws.Cells.Item(132, 11).Value = 38469183
$ws.Cells.Item(132, 12).Value = 13236
$ws.Cells.Item(132, 13).Value = -38466653
$ws.Cells.Item(132, 14).Value = -18296
$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 10).Value = 0
$ws.Cells.Item(133, 12).Value = 0
$ws.Cells.Item(133, 14).ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 5322.7393
$ws.Cells.Item(122, 9).Value = 4776.125
$ws.Cells.Item(122, 11).Value = 14328.375
$ws.Cells.Item(122, 13).Value = -11878.375
$ws.Cells.Item(133, 8).Value = 89000
$ws.Cells.Item(133, 10).Value = 89000
$ws.Cells.Item(133, 12).Value = 89000
$ws.Cells.Item(133, 14).Value = -94060
$ws.Cells.Item(136, 8).Value = 6199
$ws.Cells.Item(136, 9).Value = 4567.6665
$ws.Cells.Item(136, 10).Value = 8156.6
$ws.Cells.Item(136, 11).Value = 13702.9995
$ws.Cells.Item(136, 12).Value = 24469.8
$ws.Cells.Item(136, 13).Value = -11152.9995
$ws.Cells.Item(136, 14).Value = -29569.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 3437.125
$ws.Cells.Item(81, 9).Value = 1839.8
$ws.Cells.Item(81, 10).Value = 6099.3335
$ws.Cells.Item(81, 11).Value = 3679.6
$ws.Cells.Item(81, 12).Value = 12198.667
$ws.Cells.Item(81, 13).Value = -2618.6
$ws.Cells.Item(81, 14).Value = -14320.667
$ws.Cells.Item(84, 8).Value = 3437.125
$ws.Cells.Item(84, 9).Value = 1839.8
$ws.Cells.Item(84, 10).Value = 6099.3335
$ws.Cells.Item(84, 11).Value = 18398
$ws.Cells.Item(84, 12).Value = 60993.335
$ws.Cells.Item(84, 13).Value = -13094
$ws.Cells.Item(84, 14).Value = -71601.33499999999
$ws.Cells.Item(97, 8).Value = 39999
$ws.Cells.Item(97, 10).Value = 39999
$ws.Cells.Item(97, 12).Value = 39999
$ws.Cells.Item(97, 14).Value = -41981
$ws.Cells.Item(113, 8).Value = 425.5862
$ws.Cells.Item(113, 9).Value = 315.33334
$ws.Cells.Item(113, 10).Value = 606
$ws.Cells.Item(113, 11).Value = 946.0000200000001
$ws.Cells.Item(113, 12).Value = 1818
$ws.Cells.Item(113, 13).Value = 1223.99998
$ws.Cells.Item(113, 14).Value = -6158
$ws.Cells.Item(132, 8).Value = 1046.341
$ws.Cells.Item(132, 9).Value = 676.8
$ws.Cells.Item(132, 11).Value = 2030.4
$ws.Cells.Item(132, 13).Value = 499.6000000000001
$ws.Cells.Item(136, 8).Value = 22223770
$ws.Cells.Item(136, 9).Value = 34723300
$ws.Cells.Item(136, 11).Value = 104169900
$ws.Cells.Item(136, 13).Value = -104167350
